$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.895.42"
$ws.Range("E2").Value = "  -2.32%  "
$ws.Range("D3").Value = "3.759.21"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'403.57"
$ws.Range("E5").Value = "  -4.81%  "
$ws.Range("D6").Value = "'131.74"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").Value = "3.746.26"
$ws.Range("E7").Value = "  +0.85%  "
$ws.Range("E8").Value = "  -6.38%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").Value = "'0.717"
$ws.Range("E10").Value = "  -6.77%  "
$ws.Range("D11").Value = "'0.168"
$ws.Range("E11").Value = "  -9.74%  "
$ws.Range("D12").Value = "'0.0000357"
$ws.Range("E12").Value = "  -12.14%  "
$ws.Range("D13").Value = "'40.61"
$ws.Range("E13").Value = "  -6.08%  "
$ws.Range("D14").Value = "4.367.73"
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("E15").Value = "  -5.23%  "
$ws.Range("D16").Value = "'14.68"
$ws.Range("E16").Value = "  +12.27%  "
$ws.Range("E17").Value = "  -1.55%  "
$ws.Range("D18").Value = "3.763.50"
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("D19").Value = "'19.41"
$ws.Range("E19").Value = "  -6.97%  "
$ws.Range("D20").Value = "66.106.06"
$ws.Range("E20").Value = "  -1.98%  "
$ws.Range("E21").Value = "  -6.48%  "
$ws.Range("D22").Value = "'409.04"
$ws.Range("E22").Value = "  -9.75%  "
$ws.Range("D23").Value = "'14.35"
$ws.Range("E23").Value = "  -8.11%  "
$ws.Range("D24").Value = "'84.70"
$ws.Range("E24").Value = "  -5.87%  "
$ws.Range("E25").Value = "  -4.99%  "
$ws.Range("D26").Value = "'5.73"
$ws.Range("E26").Value = "  +14.55%  "
$ws.Range("D27").Value = "'35.92"
$ws.Range("E27").Value = "  -6.06%  "
$ws.Range("E28").Value = "  -7.26%  "
$ws.Range("D29").Value = "'9.36"
$ws.Range("E29").Value = "  -8.40%  "
$ws.Range("D30").Value = "'12.30"
$ws.Range("E30").Value = "  -2.80%  "
$ws.Range("E31").Value = "  -0.33%  "
$ws.Range("E32").Value = "  -4.74%  "
$ws.Range("D33").Value = "'7.36"
$ws.Range("E33").Value = "  -1.24%  "
$ws.Range("E34").Value = "  -6.69%  "
$ws.Range("D35").Value = "'38.90"
$ws.Range("E35").Value = "  -7.84%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").Value = "'54.74"
$ws.Range("E37").Value = "  -3.25%  "
$ws.Range("D38").Value = "0.0₃0737"
$ws.Range("E38").Value = "  -7.39%  "
$ws.Range("E39").Value = "  -7.70%  "
$ws.Range("D40").Value = "'2.88"
$ws.Range("E40").Value = "  -7.38%  "
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("E42").Value = "  -8.77%  "
$ws.Range("D43").Value = "'27.07"
$ws.Range("E43").Value = "  -1.64%  "
$ws.Range("D44").Value = "'144.76"
$ws.Range("E44").Value = "  -2.23%  "
$ws.Range("D45").Value = "'3.12"
$ws.Range("E45").Value = "  +18.32%  "
$ws.Range("D46").Value = "'3.23"
$ws.Range("E46").Value = "  -5.60%  "
$ws.Range("D47").Value = "'2.05"
$ws.Range("E47").Value = "  -2.85%  "
$ws.Range("D48").Value = "'2.58"
$ws.Range("E48").Value = "  -4.13%  "
$ws.Range("E49").Value = "  -5.39%  "
$ws.Range("D50").Value = "'2.77"
$ws.Range("E50").Value = "  -5.71%  "
$ws.Range("D51").Value = "'0.292"
$ws.Range("E51").Value = "  -5.45%  "
